$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B44").Value = 44034
$ws.Range("B44").NumberFormat = "YYYY-MM-DD"
$ws.Range("C44").Value = 16322
$ws.Range("D44").Value = 131
$ws.Range("E44").Value = 208
$ws.Range("F44").Value = 2
$ws.Range("G44").Value = 1.27
$ws.Range("H44").Value = 1.53
$ws.Range("J44").Value = $true
$ws.Range("O44").Value = "Success!"
